# Applies the "Added VWA charts and columns" edit to the Yearly sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column B: update VWA figures and apply plain "0" number format (no border) ---
$ws.Range("B2").Value = 2788.2166666666667
$ws.Range("B3").Value = 1868.9008941877794
$ws.Range("B4").Value = 2079.0556338028168
$ws.Range("B5").Value = 2408.5410000000002
$ws.Range("B6").Value = 2531.4380768461538
$ws.Range("B2:B6").NumberFormat = "0"
$ws.Range("B2:B6").Borders.LineStyle = 0

# --- New "installed power" header (new shared string) ---
$ws.Range("H1").Value = "installed power"
$ws.Range("H1").Style = "Normal"

# --- New G column: plain Year helper values (no special style) ---
$ws.Range("G2").Value = 2020
$ws.Range("G3").Value = 2021
$ws.Range("G4").Value = 2022
$ws.Range("G5").Value = 2023
$ws.Range("G6").Value = 2024
$ws.Range("G2:G6").Style = "Normal"

# --- New H column: installed power values ---
$ws.Range("H2").Value = 540
$ws.Range("H3").Value = 671
$ws.Range("H4").Value = 710
$ws.Range("H5").Value = 1000
$ws.Range("H6").Value = 1300
$ws.Range("H2:H4").Style = "Normal"
$ws.Range("H5:H6").NumberFormat = "#,##0"

# --- F6:F7 blank helper cells with thousands-separator format ---
$ws.Range("F6:F7").NumberFormat = "#,##0"

# --- Extend the "0"-formatted helper column further down (skip row 7) ---
$ws.Range("B8:B12").NumberFormat = "0"

# --- Sheet view: zoom to 125% and move selection to B6 ---
$ws.Range("B6").Select()
$excel.ActiveWindow.Zoom = 125
